$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Through 2022-05-02")

# Rename the sheet (also updates the TitlesOfParts in app.xml)
$ws.Name = "Through 2022-05-03"

# Update the header label in I1
$ws.Range("I1").Value = "2022 (through 05-03)"

# Update the two changed data values
$ws.Range("I6").Value = 10
$ws.Range("I14").Value = 561
